$wb = $excel.ActiveWorkbook

# --- Update the "conversion" note text on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$old = $ws1.Range("A1").Value()
$new = $old.Replace("✅ 1000 Bs = 9.85 = 41822.66 pesos", "✅ 1000 Bs = 9.79 = 41136.14 pesos")
$new = $new.Replace("✅ 41822.66 pesos = 9.89 = 979.32 Bs", "✅ 41136.14 pesos = 9.72 = 962.77 Bs")
$ws1.Range("A1").Value = $new

# --- Update the rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 102.1
$ws2.Range("O10").Value = 4200
$ws2.Range("N12").Value = 4229.97
$ws2.Range("O12").Value = 99
